$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the comment on AE69 (row 69 had an empty AE69 cell just holding a comment)
$ws.Range("AE69").Comment.Delete()

# 2) Rename retailer alias "Ranveer" -> "Ranveer/Rina" (cell B17)
$ws.Range("B17").Value = "Ranveer/Rina"

# 3) Widen column F slightly
$ws.Columns.Item(6).ColumnWidth = 8.57142857142857

# 4) Update the active selection (bottom-right frozen pane) to D10
$ws.Application.ActiveWindow.ScrollColumn = 22
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("D10").Select()

# 5) Fill in newly-collected data for 23 Jan (column AC) across many rows
$acValues = @{
    7  = 8000
    15 = 2000
    17 = 3000
    19 = 2000
    23 = 3000
    25 = 1000
    31 = 500
    38 = 2000
    46 = 2300
    48 = 3000
    49 = 3000
    52 = 1000
    53 = 5000
    55 = 5500
    61 = 4000
    63 = 2000
    65 = 5000
    67 = 2000
    68 = 1200
    70 = 3000
    77 = 2000
    78 = 2000
    82 = 5000
}

# Rows whose new AC cell uses the highlighted ("Cash") fill style, matched from
# the existing highlighted cell H3 used elsewhere in the sheet for that style.
$highlightRows = @(7, 17, 19, 38, 49, 77, 82)

foreach ($row in $acValues.Keys) {
    $cell = $ws.Cells.Item($row, 29)
    $cell.Value = $acValues[$row]
    if ($highlightRows -contains $row) {
        $cell.Interior.Color = $ws.Range("H3").Interior.Color
    }
}

$wb.Application.Calculate()
